# cambios de las fracciones
# Update the quarterly reporting dates (3er Trimestre -> 4to Trimestre) on
# the "Reporte de Formatos" sheet, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Row 8 period/update dates: shift from Q3 2022 to Q4 2022
$ws.Range("B8").Value = Get-Date -Year 2022 -Month 10 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("C8").Value = Get-Date -Year 2022 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0
$ws.Range("I8").Value = Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("J8").Value = Get-Date -Year 2023 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0

# Move the active cell selection from C14 to C11
$ws.Activate()
$ws.Range("C11").Select()
